$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Header row (row 1): turn the old duplicate-data "header" into real column
# headers, and extend it with the same metadata columns used on the other
# asset sheets (stock, etc.) ---
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

# Copy the header style (bold + border) onto the newly added header cells.
$ws.Range("F1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)

$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# --- Data rows: B:F already hold the correct values; append the metadata
# columns G:M (mirroring sheet4/"股票"'s layout) ---

# Copy the plain data-row style onto the new cells for rows 2-4 first.
$ws.Range("F2").Copy()
$ws.Range("G2:M4").PasteSpecial(-4122)

# Row 2
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "'2011-11-22"
$ws.Range("J2").Value = "馬文君"
$ws.Range("K2").Value = 1724
$ws.Range("L2").Value = "tmp99351"
$ws.Range("M2").Value = 55

# Row 3
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "'2011-11-22"
$ws.Range("J3").Value = "馬文君"
$ws.Range("K3").Value = 1724
$ws.Range("L3").Value = "tmp99351"
$ws.Range("M3").Value = 56

# Row 4
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "'2011-11-22"
$ws.Range("J4").Value = "馬文君"
$ws.Range("K4").Value = 1724
$ws.Range("L4").Value = "tmp99351"
$ws.Range("M4").Value = 57
